# Major accuracy check update:
# 1. Update the poly(A) isolation protocol value for all sample rows
#    from "NEBNextPoly(A)E7490" to "NEBNextPoly(A)E7490L" (column G, rows 2-13).
# 2. Widen column G to fit the longer value.
# 3. Move the active selection from column I to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the cell values in column G (rows 2 through 13).
$ws.Range("G2:G13").Value = "NEBNextPoly(A)E7490L"

# 2. Widen column G to accommodate the updated text.
$ws.Range("G1").ColumnWidth = 27.78

# 3. Update the selected range/active cell to column G.
$ws.Range("G2:G13").Select() | Out-Null
